$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Title heading and the later bold repetition of the same title (both occurrences)
Replace-Text "Play Great 88 for Free: A Bonus-Rich Slot Game" "Play Great 88 Slot Game for Free"

# "What we like" bullet list
Replace-Text "Bonus-rich game with several opportunities to win big" "Impressive graphics and immersive oriental music"
Replace-Text "Stunning graphics and immersive oriental music" "Bonus-rich game with multiple opportunities to win big"
Replace-Text "Wide range of betting options for all types of players" "Wide range of betting options"
Replace-Text "Potential payouts of up to 5,000 times your bet" "Autoplay function and gamble feature for added excitement"

# "What we don't like" bullet list
Replace-Text "RTP is average compared to other slot games" "RTP is average, not standout"
Replace-Text "Coins values are relatively low" "Limited paylines compared to some other slot games"

# Closing italic summary paragraph
Replace-Text "Experience traditional Chinese celebrations with immersive graphics and audio in Great 88. Play for free and try to win big with several bonus opportunities." "Discover the traditional Chinese celebrations in Great 88 slot game. Play for free and win big!"
